# Gym Management System ER Table - remove User_ID column from the
# Email_IDs sheet (rename to Email_ID), and update the various sheet
# selections / active sheet to match the author's final view state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Member" sheet: selection changes to a single cell B14
# ---------------------------------------------------------------------
$wsMember = $wb.Worksheets.Item("Member")
$null = $wsMember.Range("B14").Select()

# ---------------------------------------------------------------------
# 2. "Personal_Instructor" sheet: selection changes to a single cell F12
#    (it also loses tabSelected, which happens automatically once a
#    different sheet becomes active further below)
# ---------------------------------------------------------------------
$wsPI = $wb.Worksheets.Item("Personal_Instructor")
$null = $wsPI.Range("F12").Select()

# ---------------------------------------------------------------------
# 3. "Email_IDs" sheet: drop the User_ID column (column A), shifting
#    Ins_ID (old column B) into column A and Email_ID (old column C)
#    into column B, then clear out the now-unused column C.
# ---------------------------------------------------------------------
$wsEmail = $wb.Worksheets.Item("Email_IDs")

for ($r = 1; $r -le 23; $r++) {
    $insIdVal = $wsEmail.Cells.Item($r, 2).Value2
    $emailVal = $wsEmail.Cells.Item($r, 3).Value2
    $wsEmail.Cells.Item($r, 1).Value = $insIdVal
    $wsEmail.Cells.Item($r, 2).Value = $emailVal
}
$wsEmail.Range("C1:C23").Clear()

# With the User_ID column gone the remaining two columns are wider
# relative to the wrapped text, so the row heights for the instructor
# rows grow (most become 43.2pt, the two "Everest" rows become 57.6pt).
for ($r = 12; $r -le 23; $r++) {
    if ($r -eq 15 -or $r -eq 16) {
        $wsEmail.Rows.Item($r).RowHeight = 57.6
    } else {
        $wsEmail.Rows.Item($r).RowHeight = 43.2
    }
}

# Rename the sheet and make it the active / selected tab with the new
# selection covering the remaining Ins_ID column.
$wsEmail.Name = "Email_ID"
$wsEmail.Activate()
$null = $wsEmail.Range("A1:A23").Select()
